# Effects sheet ("xl/worksheets/sheet5.xml") gains three new blocks of
# Enabled/Value field pairs (STR..WIS, CHR/LUCK, hp/totalHP/mana,
# earth/fire/water/lightning DR) in the middle of the existing table, plus a
# header-row border/alignment fix. Data-driven: build the full target table
# (row, A-label, A-style-class, B-value, B-style-class, C-value, C-style-class)
# and replay it onto the sheet after making room for the new rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Effects")

# --- 1. Make room: insert 30 blank rows before the old row 13 -------------
#     Old layout: 1 header, 2 name, 3-12 damage/range/target fields,
#                 13-24 ac/damageMod/bluntDR/chopDR/pierceDR/slashDR pairs,
#                 25 description.
#     New layout adds a STR..WIS/CHR/LUCK block (16 rows) + a
#     hp/totalHP/mana block (6 rows) + an earth/fire/water/lightning DR
#     block (8 rows) = 30 new rows, pushing the old 13-25 block down to
#     29-55 (in three pieces, since the new blocks are interleaved).
$ws.Range("A13:A42").EntireRow.Insert()

# --- 2. Full target content ------------------------------------------------
# columns: row, A text, A style-class, B value, B style-class, C value, C style-class
# style-class 1 = bordered (left-aligned), 2/3 = bold+fill+bordered label cell,
# 5 = bordered + right-aligned
$rowsData = @(
    @(1, "ID", 2, 1, 5, 2, 5),
    @(2, "name", 2, "Fireball", 5, "Burn", 5),
    @(3, "damageTypeEnabled", 2, 1, 5, 1, 1),
    @(4, "damageType", 3, "f", 5, "f", 5),
    @(5, "rangeEnabled", 3, 1, 5, 0, 1),
    @(6, "range", 3, "1,0,10,1,20,2,30,3", 5, 0, 1),
    @(7, "targetedEnabled", 3, 1, 5, 0, 1),
    @(8, "targeted", 3, "0,-1,1,0", 5, 0, 1),
    @(9, "diceDamageEnabled", 3, 1, 5, 1, 1),
    @(10, "diceDamage", 3, "4,1,6,2,8,3,10,4,12,5", 5, "4,1,6,2,8,3,10,4,12,5", 5),
    @(11, "diceDamagedurationEnabled", 3, 0, 5, 1, 1),
    @(12, "diceDamageDuration", 3, 0, 5, "4,1,6,2,8,3,10,4,12,5", 5),
    @(13, "STREnabled", 3, 0, 5, 0, 5),
    @(14, "STR", 3, 0, 5, 0, 5),
    @(15, "DEXEnabled", 3, 0, 5, 0, 5),
    @(16, "DEX", 3, 0, 5, 0, 5),
    @(17, "CONEnabled", 3, 0, 5, 0, 5),
    @(18, "CON", 3, 0, 5, 0, 5),
    @(19, "WILLEnabled", 3, 0, 5, 0, 5),
    @(20, "WILL", 3, 0, 5, 0, 5),
    @(21, "INTEnabled", 3, 0, 5, 0, 5),
    @(22, "INT", 3, 0, 5, 0, 5),
    @(23, "WISEnabled", 3, 0, 5, 0, 5),
    @(24, "WIS", 3, 0, 5, 0, 5),
    @(25, "CHREnabled", 3, 0, 5, 0, 5),
    @(26, "CHR", 3, 0, 5, 0, 5),
    @(27, "LUCKEnabled", 3, 0, 5, 0, 5),
    @(28, "LUCK", 3, 0, 5, 0, 5),
    @(29, "acEnabled", 3, 0, 5, 0, 1),
    @(30, "ac", 3, 0, 5, 0, 5),
    @(31, "damageModEnabled", 3, 1, 5, 0, 1),
    @(32, "damageMod", 3, "0,0,1,1,2,2,3,3,4,4,5,5,6,6,7,7,8,8,9,9,10,10", 5, 0, 5),
    @(33, "hpEnabled", 3, 0, 5, 0, 5),
    @(34, "hp", 3, 0, 5, 0, 5),
    @(35, "totalHPEnabled", 3, 0, 5, 0, 5),
    @(36, "totalHP", 3, 0, 5, 0, 5),
    @(37, "manaEnabled", 3, 0, 5, 0, 5),
    @(38, "mana", 3, 0, 5, 0, 5),
    @(39, "bluntDREnabled", 3, 0, 5, 0, 1),
    @(40, "bluntDR", 3, 0, 5, 0, 5),
    @(41, "chopDREnabled", 3, 0, 5, 0, 1),
    @(42, "chopDR", 3, 0, 5, 0, 5),
    @(43, "pierceDREnabled", 3, 0, 5, 0, 1),
    @(44, "pierceDR", 3, 0, 5, 0, 5),
    @(45, "slashDREnabled", 3, 0, 5, 0, 1),
    @(46, "slashDR", 3, 0, 5, 0, 5),
    @(47, "earthDREnabled", 3, 0, 5, 0, 5),
    @(48, "earthDR", 3, 0, 5, 0, 5),
    @(49, "fireDREnabled", 3, 0, 5, 0, 5),
    @(50, "fireDR", 3, 0, 5, 0, 5),
    @(51, "waterDREnabled", 3, 0, 5, 0, 5),
    @(52, "waterDR", 3, 0, 5, 0, 5),
    @(53, "lightningDREnabled", 3, 0, 5, 0, 5),
    @(54, "lightningDR", 3, 0, 5, 0, 5),
    @(55, "description", 3, "Targeted effect: a flying ball of fire", 5, "Duration effect: fire damage over time", 1)
)

# --- 3. Style donors: cells whose formatting never changes in this edit ---
$styleDonorA = $ws.Cells.Item(4, 1)   # style "3" - bold/filled/bordered label
$styleDonorB = $ws.Cells.Item(4, 2)   # style "5" - bordered + right aligned
$styleDonorC1 = $ws.Cells.Item(5, 3)  # style "1" - bordered only

function Copy-Style($srcCell, $dstCell) {
    $srcCell.Copy() | Out-Null
    $dstCell.PasteSpecial(-4122) | Out-Null
}

function Apply-StyleClass($cell, [int]$styleClass) {
    if ($styleClass -eq 1) {
        Copy-Style $styleDonorC1 $cell
    } elseif ($styleClass -eq 2 -or $styleClass -eq 3) {
        Copy-Style $styleDonorA $cell
    } elseif ($styleClass -eq 5) {
        Copy-Style $styleDonorB $cell
    }
}

foreach ($row in $rowsData) {
    $r = $row[0]
    $aVal = $row[1]
    $aStyle = $row[2]
    $bVal = $row[3]
    $bStyle = $row[4]
    $cVal = $row[5]
    $cStyle = $row[6]

    $aCell = $ws.Cells.Item($r, 1)
    $aCell.Value = $aVal
    Apply-StyleClass $aCell $aStyle

    $bCell = $ws.Cells.Item($r, 2)
    $bCell.Value = $bVal
    Apply-StyleClass $bCell $bStyle

    $cCell = $ws.Cells.Item($r, 3)
    $cCell.Value = $cVal
    Apply-StyleClass $cCell $cStyle
}

$excel.CutCopyMode = 0

# --- 4. Sheet view: scroll position + selection ----------------------------
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 23
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A31").Select()
